$d = $word.ActiveDocument

# The document originally has a single paragraph containing "dfghhgdfh"
# (plus a _GoBack bookmark). The edit:
#   1) changes that paragraph's text to "gfbnfgbnfgnfnf"
#   2) adds a new paragraph right after it, with the original text
#      "dfghhgdfh" (same en-US language run formatting).

$p1 = $d.Paragraphs(1)

# Insert a new paragraph after the first one, then give it the
# original "dfghhgdfh" text.
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "dfghhgdfh"

# Now update the first paragraph's own text, scoping Find to its range
# only (so it doesn't also match the text we just inserted below it).
$p1.Range.Find.Execute("dfghhgdfh", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "gfbnfgbnfgnfnf", 2)
